$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number formats/styles into the new D:E columns from the (now shifted) F:G columns
$ws.Range("F7:F35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("G7:G35").Copy()
$ws.Range("E7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("G38:G77").Copy()
$ws.Range("E38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Range("G80:G102").Copy()
$ws.Range("E80:E102").PasteSpecial(-4122)

# Populate the new quarter data in columns D and E
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 141500
$ws.Range("E8").Value = 149300
$ws.Range("D9").Value = 104000
$ws.Range("E9").Value = 109000
$ws.Range("D10").Value = 37500
$ws.Range("E10").Value = 40300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 1800
$ws.Range("E14").Value = 200
$ws.Range("D15").Value = 23000
$ws.Range("E15").Value = 23500
$ws.Range("D17").Value = 145300
$ws.Range("E17").Value = 145000
$ws.Range("D18").Value = -3800
$ws.Range("E18").Value = 4300
$ws.Range("D20").Value = -300
$ws.Range("E20").Value = 500
$ws.Range("D21").Value = 19000
$ws.Range("E21").Value = 28300
$ws.Range("D22").Value = 9800
$ws.Range("E22").Value = 9800
$ws.Range("D23").Value = -13900
$ws.Range("E23").Value = -5000
$ws.Range("D24").Value = 600
$ws.Range("E24").Value = 300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -14500
$ws.Range("E26").Value = -5200
$ws.Range("D27").Value = -14500
$ws.Range("E27").Value = -5200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 300
$ws.Range("E32").Value = -500
$ws.Range("D33").Value = -14500
$ws.Range("E33").Value = -5200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -14500
$ws.Range("E35").Value = -5200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 53600
$ws.Range("E41").Value = 51500
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 130900
$ws.Range("E43").Value = 139700
$ws.Range("D44").Value = 18900
$ws.Range("E44").Value = 19000
$ws.Range("D45").Value = 11700
$ws.Range("E45").Value = 13700
$ws.Range("D46").Value = 215000
$ws.Range("E46").Value = 223900
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 524900
$ws.Range("E48").Value = 527300
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1700
$ws.Range("E52").Value = 1700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 741600
$ws.Range("E54").Value = 752900
$ws.Range("D57").Value = 34100
$ws.Range("E57").Value = 34700
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 70600
$ws.Range("E59").Value = 69100
$ws.Range("D60").Value = 104800
$ws.Range("E60").Value = 103800
$ws.Range("D61").Value = 464600
$ws.Range("E61").Value = 463800
$ws.Range("D62").Value = 7200
$ws.Range("E62").Value = 6700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 576500
$ws.Range("E66").Value = 574400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -388400
$ws.Range("E72").Value = -373900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 165100
$ws.Range("E76").Value = 178500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -14500
$ws.Range("E81").Value = -5200
$ws.Range("D83").Value = 23000
$ws.Range("E83").Value = 23500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 18200
$ws.Range("E89").Value = 4400
$ws.Range("D91").Value = -18400
$ws.Range("E91").Value = -17300
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -17100
$ws.Range("E94").Value = -14400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 1100
$ws.Range("E102").Value = -10000

# Row 49 (Other Assets, Current) had several previously-"NA" quarters newly reported as 0
# after the refresh; columns F:J (old D:H, shifted right by 2) need updating from "NA" to 0.
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
